$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Run 50" column (column AZ). This shifts the old "Mean"
# column (BA) left into AZ, and shifts nothing else since it's the
# last-but-one column.
$ws.Columns.Item(52).Delete()

# Header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A (now the MaxFES fractions instead of generation counts)
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Recomputed "Mean" column now lives in AZ (after the column delete above)
$ws.Range("AZ2").Value = 353233832.9672861
$ws.Range("AZ3").Value = 243545136.913418
$ws.Range("AZ4").Value = 46022831.13770574
$ws.Range("AZ5").Value = 3091685.76278633
$ws.Range("AZ6").Value = 1148607.90251364
$ws.Range("AZ7").Value = 700078.28426107
$ws.Range("AZ8").Value = 476350.60460794
$ws.Range("AZ9").Value = 328402.8122585
$ws.Range("AZ10").Value = 230276.52377992
$ws.Range("AZ11").Value = 185154.24899913
$ws.Range("AZ12").Value = 144794.32127092
$ws.Range("AZ13").Value = 120952.4720088
$ws.Range("AZ14").Value = 92477.79611323
